$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Row 3 corresponds to the "child_id" variable.
# Change valueType (column B) from "integer" to "text"
$ws.Range("B3").Value = "text"

# Clear the unit (column C) value - child_id no longer has a unit ("numeric")
$ws.Range("C3").ClearContents()
